# "Proyecto completo" — move the Curve-Fitting data table from A1:B21 down
# and to the right (to C4:D24), box it with a thin border, and re-fit the
# column widths, matching the layout tweak the author made before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the table from A1:B21 to C4:D24 ------------------------------
$ws.Range("A1:B21").Copy() | Out-Null
$ws.Range("C4").PasteSpecial() | Out-Null
$ws.Range("A1:B21").ClearContents() | Out-Null
$excel.CutCopyMode = 0

# --- Box the moved table with a thin border on all four sides ----------
$tableRange = $ws.Range("C4:D24")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# --- Re-fit column widths (old A/B columns keep their fitted widths, ---
# --- new C/D columns pick up slightly wider ones once bordered) --------
$ws.Columns.Item(1).ColumnWidth = 11.92
$ws.Columns.Item(2).ColumnWidth = 22.59
$ws.Columns.Item(3).ColumnWidth = 12.26
$ws.Columns.Item(4).ColumnWidth = 22.92

# --- Match the final on-screen selection --------------------------------
$ws.Range("F23").Select() | Out-Null
